$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: H6, I6, J6, K6, L6, M6, N6
$ws.Range("H6").Value = 349.1
$ws.Range("I6").Value = 254.55556
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 763.66668
$ws.Range("L6").Value = 3600
$ws.Range("M6").Value = -651.66668
$ws.Range("N6").Value = -3824

# Row 69: H69, I69, K69, M69
$ws.Range("H69").Value = 4876.923
$ws.Range("I69").Value = 4822.222
$ws.Range("K69").Value = 14466.666
$ws.Range("M69").Value = -13592.666

# Row 72: H72, I72, K72, M72
$ws.Range("H72").Value = 4876.923
$ws.Range("I72").Value = 4822.222
$ws.Range("K72").Value = 43399.998
$ws.Range("M72").Value = -39031.998

# Row 76: H76, I76, J76, K76, L76, M76, N76
$ws.Range("H76").Value = 6000.4287
$ws.Range("I76").Value = 6333.8335
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 6333.8335
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -6018.8335
$ws.Range("N76").Value = -4630

# Row 79: H79, I79, J79, K79, L79, M79, N79
$ws.Range("H79").Value = 6000.4287
$ws.Range("I79").Value = 6333.8335
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 6333.8335
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -5241.8335
$ws.Range("N79").Value = -6184

# Row 88: H88, J88, L88, N88
$ws.Range("H88").Value = 3690
$ws.Range("J88").Value = 3690
$ws.Range("L88").Value = 3690
$ws.Range("N88").Value = -4502

# Row 91: H91, J91, L91, N91
$ws.Range("H91").Value = 3690
$ws.Range("J91").Value = 3690
$ws.Range("L91").Value = 3690
$ws.Range("N91").Value = -6498

# Row 131: H131, I131, J131, K131, L131, M131, N131
$ws.Range("H131").Value = 6508.4595
$ws.Range("I131").Value = 1217.5385
$ws.Range("J131").Value = 9374.375
$ws.Range("K131").Value = 3652.6155
$ws.Range("L131").Value = 28123.125
$ws.Range("M131").Value = 1387.3845
$ws.Range("N131").Value = -38203.125

# Row 138: H138, I138, J138, K138, L138, M138, N138
$ws.Range("H138").Value = 3815.3948
$ws.Range("I138").Value = 2253.25
$ws.Range("J138").Value = 4536.385
$ws.Range("K138").Value = 6759.75
$ws.Range("L138").Value = 13609.155
$ws.Range("M138").Value = -1619.75
$ws.Range("N138").Value = -23889.155

$ws = $wb.Worksheets.Item("ARM")
# Row 37: H37, J37, L37 | clear: N37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

# Row 122: H122, I122, J122, K122, L122, M122 | clear: N122
$ws.Range("H122").Value = 2141.75
$ws.Range("I122").Value = 2141.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6425.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3975.25
$ws.Range("N122").ClearContents()

# Row 131: H131, J131, L131, N131
$ws.Range("H131").Value = 59713.855
$ws.Range("J131").Value = 59713.855
$ws.Range("L131").Value = 59713.855
$ws.Range("N131").Value = -69793.85500000001

# Row 140: H140, I140, J140, K140, L140, N140 | clear: M140
$ws.Range("H140").Value = 66500
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 66500
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 66500
$ws.Range("N140").Value = -76860
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 11: H11, I11, J11, K11, L11, M11 | clear: N11
$ws.Range("H11").Value = 552
$ws.Range("I11").Value = 552
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 552
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -412
$ws.Range("N11").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 11: H11, I11, K11, M11
$ws.Range("H11").Value = 49999.5
$ws.Range("I11").Value = 49999.5
$ws.Range("K11").Value = 49999.5
$ws.Range("M11").Value = -49859.5

# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 2639.1365
$ws.Range("I132").Value = 2038.9412
$ws.Range("K132").Value = 6116.8236
$ws.Range("M132").Value = -3586.8236

$ws = $wb.Worksheets.Item("CUL")
# Row 3: H3, I3, J3, K3, L3, M3, N3
$ws.Range("H3").Value = 5947.2
$ws.Range("I3").Value = 2909.75
$ws.Range("J3").Value = 7972.1665
$ws.Range("K3").Value = 8729.25
$ws.Range("L3").Value = 23916.4995
$ws.Range("M3").Value = -8617.25
$ws.Range("N3").Value = -24140.4995

# Row 4: H4, I4, J4, K4, L4, M4 | clear: N4
$ws.Range("H4").Value = 350050
$ws.Range("I4").Value = 350050
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1050150
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1050038
$ws.Range("N4").ClearContents()

# Row 103: H103, J103, L103, N103
$ws.Range("H103").Value = 2783.875
$ws.Range("J103").Value = 3130.4614
$ws.Range("L103").Value = 9391.3842
$ws.Range("N103").Value = -11149.3842

# Row 122: H122, J122, L122, N122
$ws.Range("H122").Value = 1160
$ws.Range("J122").Value = 1225
$ws.Range("L122").Value = 11025
$ws.Range("N122").Value = -15925

$ws = $wb.Worksheets.Item("GSM")
# Row 10: H10, I10, K10, M10
$ws.Range("H10").Value = 17163.334
$ws.Range("I10").Value = 845
$ws.Range("K10").Value = 845
$ws.Range("M10").Value = -676

# Row 80: H80, I80, J80, K80, L80, M80, N80
$ws.Range("H80").Value = 3316.2942
$ws.Range("I80").Value = 3489.5557
$ws.Range("J80").Value = 3121.375
$ws.Range("K80").Value = 3489.5557
$ws.Range("L80").Value = 3121.375
$ws.Range("M80").Value = -2491.5557
$ws.Range("N80").Value = -5117.375

# Row 83: H83, I83, J83, K83, L83, M83, N83
$ws.Range("H83").Value = 3316.2942
$ws.Range("I83").Value = 3489.5557
$ws.Range("J83").Value = 3121.375
$ws.Range("K83").Value = 17447.7785
$ws.Range("L83").Value = 15606.875
$ws.Range("M83").Value = -12455.7785
$ws.Range("N83").Value = -25590.875

$ws = $wb.Worksheets.Item("LTW")
# Row 7: H7, I7, J7, K7, L7, M7, N7
$ws.Range("H7").Value = 4782.778
$ws.Range("I7").Value = 4173.3335
$ws.Range("J7").Value = 6001.6665
$ws.Range("K7").Value = 4173.3335
$ws.Range("L7").Value = 6001.6665
$ws.Range("M7").Value = -4061.3335
$ws.Range("N7").Value = -6225.6665

# Row 46: H46, I46, J46, K46, L46, M46, N46
$ws.Range("H46").Value = 1635.2727
$ws.Range("I46").Value = 1396
$ws.Range("J46").Value = 1725
$ws.Range("K46").Value = 1396
$ws.Range("L46").Value = 1725
$ws.Range("M46").Value = -1208
$ws.Range("N46").Value = -2101

# Row 55: H55, I55, J55, K55, L55, M55, N55
$ws.Range("H55").Value = 515
$ws.Range("I55").Value = 150
$ws.Range("J55").Value = 697.5
$ws.Range("K55").Value = 150
$ws.Range("L55").Value = 697.5
$ws.Range("M55").Value = 23
$ws.Range("N55").Value = -1043.5

# Row 126: H126, I126, J126, K126, L126, M126, N126
$ws.Range("H126").Value = 4782.778
$ws.Range("I126").Value = 4173.3335
$ws.Range("J126").Value = 6001.6665
$ws.Range("K126").Value = 12520.0005
$ws.Range("L126").Value = 18004.9995
$ws.Range("M126").Value = -10050.0005
$ws.Range("N126").Value = -22944.9995

# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 4139.8438
$ws.Range("I132").Value = 4171.4346
$ws.Range("K132").Value = 12514.3038
$ws.Range("M132").Value = -9984.303799999998

$ws = $wb.Worksheets.Item("WVR")
# Row 107: H107, I107, K107, M107
$ws.Range("H107").Value = 478.14285
$ws.Range("I107").Value = 438.76923
$ws.Range("K107").Value = 1316.30769
$ws.Range("M107").Value = 603.6923099999999

# Row 122: H122, J122, L122, N122
$ws.Range("H122").Value = 43105960
$ws.Range("J122").Value = 4034.5454
$ws.Range("L122").Value = 12103.6362
$ws.Range("N122").Value = -17003.6362
